$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

$wsZhCn.Range("B3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("G3").Value = "2016-03-03 15:54:33"

$wsDeDe.Range("B3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("G3").Value = "2016-03-03 15:54:58"
